# Applies numeric value updates to columns H:N (market price / profit calculations)
# across all 8 sheets, per the Omega_Profits commit diff.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1681
$ws.Range("I28").Value = 438.16666
$ws.Range("J28").Value = 4166.6665
$ws.Range("K28").Value = 438.16666
$ws.Range("L28").Value = 4166.6665
$ws.Range("M28").Value = 46.83334000000002
$ws.Range("N28").Value = -5136.6665
# Row 31
$ws.Range("H31").Value = 2300.182
$ws.Range("I31").Value = 2030.2
$ws.Range("K31").Value = 6090.6
$ws.Range("M31").Value = -5860.6
# Row 39
$ws.Range("H39").Value = 5141.375
$ws.Range("I39").Value = 5022.5
$ws.Range("J39").Value = 5181
$ws.Range("K39").Value = 15067.5
$ws.Range("L39").Value = 15543
$ws.Range("M39").Value = -14771.5
$ws.Range("N39").Value = -16135
# Row 51
$ws.Range("H51").Value = 34335.875
$ws.Range("I51").Value = 8350
$ws.Range("J51").Value = 38048.145
$ws.Range("K51").Value = 8350
$ws.Range("L51").Value = 38048.145
$ws.Range("M51").Value = -7866
$ws.Range("N51").Value = -39016.145
# Row 64
$ws.Range("H64").Value = 8363
$ws.Range("I64").Value = 7399.6
$ws.Range("J64").Value = 9165.833000000001
$ws.Range("K64").Value = 7399.6
$ws.Range("L64").Value = 9165.833000000001
$ws.Range("M64").Value = -7151.6
$ws.Range("N64").Value = -9661.833000000001
# Row 67
$ws.Range("H67").Value = 8363
$ws.Range("I67").Value = 7399.6
$ws.Range("J67").Value = 9165.833000000001
$ws.Range("K67").Value = 7399.6
$ws.Range("L67").Value = 9165.833000000001
$ws.Range("M67").Value = -6541.6
$ws.Range("N67").Value = -10881.833
# Row 107
$ws.Range("H107").Value = 1965.4445
$ws.Range("I107").Value = 1915.4
$ws.Range("J107").Value = 2215.6667
$ws.Range("K107").Value = 1915.4
$ws.Range("L107").Value = 2215.6667
$ws.Range("M107").Value = 4.599999999999909
$ws.Range("N107").Value = -6055.6667
# Row 136
$ws.Range("H136").Value = 68749
$ws.Range("J136").Value = 68749
$ws.Range("L136").Value = 68749
$ws.Range("N136").Value = -78949
# Row 137
$ws.Range("H137").Value = 1967.1923
$ws.Range("I137").Value = 1651.909
$ws.Range("J137").Value = 2198.4
$ws.Range("K137").Value = 4955.727000000001
$ws.Range("L137").Value = 6595.200000000001
$ws.Range("M137").Value = -2405.727000000001
$ws.Range("N137").Value = -11695.2
# Row 141
$ws.Range("H141").Value = 3507.8572
$ws.Range("I141").Value = 3507.8572
$ws.Range("K141").Value = 10523.5716
$ws.Range("M141").Value = -5343.571599999999

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5310.636
$ws.Range("I61").Value = 4677.25
$ws.Range("K61").Value = 4677.25
$ws.Range("M61").Value = -4465.25
# Row 102
$ws.Range("H102").Value = 2112.5
$ws.Range("I102").Value = 2247.3076
$ws.Range("K102").Value = 2247.3076
$ws.Range("M102").Value = -625.3076000000001
# Row 132
$ws.Range("H132").Value = 1875.5358
$ws.Range("I132").Value = 1889.0769
$ws.Range("J132").Value = 1699.5
$ws.Range("K132").Value = 5667.2307
$ws.Range("L132").Value = 5098.5
$ws.Range("M132").Value = -3137.2307
$ws.Range("N132").Value = -10158.5
# Row 136
$ws.Range("H136").Value = 5310.636
$ws.Range("I136").Value = 4677.25
$ws.Range("K136").Value = 14031.75
$ws.Range("M136").Value = -11481.75

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 2565
$ws.Range("J64").Value = 2438
$ws.Range("L64").Value = 2438
$ws.Range("N64").Value = -2888
# Row 67
$ws.Range("H67").Value = 2565
$ws.Range("J67").Value = 2438
$ws.Range("L67").Value = 2438
$ws.Range("N67").Value = -3998
# Row 80
$ws.Range("H80").Value = 1312.8889
$ws.Range("I80").Value = 1636
$ws.Range("J80").Value = 989.7778
$ws.Range("K80").Value = 1636
$ws.Range("L80").Value = 989.7778
$ws.Range("M80").Value = -638
$ws.Range("N80").Value = -2985.7778
# Row 83
$ws.Range("H83").Value = 1312.8889
$ws.Range("I83").Value = 1636
$ws.Range("J83").Value = 989.7778
$ws.Range("K83").Value = 8180
$ws.Range("L83").Value = 4948.889
$ws.Range("M83").Value = -3188
$ws.Range("N83").Value = -14932.889
# Row 95
$ws.Range("H95").Value = 27892
$ws.Range("J95").Value = 27892
$ws.Range("L95").Value = 27892
$ws.Range("N95").Value = -33384

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 5200
$ws.Range("I14").Value = 5200
$ws.Range("K14").Value = 5200
$ws.Range("M14").Value = -5030
# Row 16
$ws.Range("H16").Value = 10340.5625
$ws.Range("I16").Value = 894.25
$ws.Range("J16").Value = 19786.875
$ws.Range("K16").Value = 894.25
$ws.Range("L16").Value = 19786.875
$ws.Range("M16").Value = -607.25
$ws.Range("N16").Value = -20360.875
# Row 19
$ws.Range("H19").Value = 550388.9
$ws.Range("I19").Value = 1100239.8
$ws.Range("K19").Value = 1100239.8
$ws.Range("M19").Value = -1100069.8
# Row 24
$ws.Range("H24").Value = 550388.9
$ws.Range("I24").Value = 1100239.8
$ws.Range("K24").Value = 1100239.8
$ws.Range("M24").Value = -1100069.8
# Row 99
$ws.Range("H99").Value = 6488.636
$ws.Range("I99").Value = 6194.75
$ws.Range("J99").Value = 6656.5713
$ws.Range("K99").Value = 6194.75
$ws.Range("L99").Value = 6656.5713
$ws.Range("M99").Value = -4696.75
$ws.Range("N99").Value = -9652.5713
# Row 113
$ws.Range("H113").Value = 10340.5625
$ws.Range("I113").Value = 894.25
$ws.Range("J113").Value = 19786.875
$ws.Range("K113").Value = 894.25
$ws.Range("L113").Value = 19786.875
$ws.Range("M113").Value = 1275.75
$ws.Range("N113").Value = -24126.875
# Row 126
$ws.Range("H126").Value = 6488.636
$ws.Range("I126").Value = 6194.75
$ws.Range("J126").Value = 6656.5713
$ws.Range("K126").Value = 18584.25
$ws.Range("L126").Value = 19969.7139
$ws.Range("M126").Value = -16114.25
$ws.Range("N126").Value = -24909.7139
# Row 132
$ws.Range("H132").Value = 3679.0576
$ws.Range("I132").Value = 3455.3777
$ws.Range("K132").Value = 10366.1331
$ws.Range("M132").Value = -7836.133099999999
# Row 141
$ws.Range("H141").Value = 62575.09
$ws.Range("J141").Value = 64703.2
$ws.Range("L141").Value = 64703.2
$ws.Range("N141").Value = -75063.2

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 72.818184
$ws.Range("I6").Value = 55.1
$ws.Range("K6").Value = 165.3
$ws.Range("M6").Value = -52.30000000000001
# Row 51
$ws.Range("H51").Value = 2944.1667
$ws.Range("I51").Value = 2944.1667
$ws.Range("K51").Value = 8832.500100000001
$ws.Range("M51").Value = -8372.500100000001
# Row 56
$ws.Range("H56").Value = 7173.273
$ws.Range("I56").Value = 7173.273
$ws.Range("K56").Value = 7173.273
$ws.Range("M56").Value = -6643.273
# Row 68
$ws.Range("H68").Value = 699
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 699
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 98
$ws.Range("H98").Value = 1512.5454
$ws.Range("I98").Value = 1185.5
$ws.Range("J98").Value = 1905
$ws.Range("K98").Value = 3556.5
$ws.Range("L98").Value = 5715
$ws.Range("M98").Value = -2058.5
$ws.Range("N98").Value = -8711
# Row 107
$ws.Range("H107").Value = 887.64703
$ws.Range("J107").Value = 1476.8889
$ws.Range("L107").Value = 4430.6667
$ws.Range("N107").Value = -8270.6667
# Row 125
$ws.Range("H125").Value = 5000
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 129
$ws.Range("H129").Value = 1071.2
$ws.Range("I129").Value = 825.36365
$ws.Range("J129").Value = 1747.25
$ws.Range("K129").Value = 2476.09095
$ws.Range("L129").Value = 5241.75
$ws.Range("M129").Value = 2523.90905
$ws.Range("N129").Value = -15241.75
# Row 131
$ws.Range("H131").Value = 1865.6123
$ws.Range("I131").Value = 777.375
$ws.Range("J131").Value = 2393.2424
$ws.Range("K131").Value = 2332.125
$ws.Range("L131").Value = 7179.7272
$ws.Range("M131").Value = 2707.875
$ws.Range("N131").Value = -17259.7272
# Row 132
$ws.Range("H132").Value = 1899.1666
$ws.Range("I132").Value = 1498.5
$ws.Range("J132").Value = 2099.5
$ws.Range("K132").Value = 13486.5
$ws.Range("L132").Value = 18895.5
$ws.Range("M132").Value = -10956.5
$ws.Range("N132").Value = -23955.5
# Row 139
$ws.Range("H139").Value = 6193.5713
$ws.Range("I139").Value = 2565.0667
$ws.Range("J139").Value = 10380.308
$ws.Range("K139").Value = 7695.2001
$ws.Range("L139").Value = 31140.924
$ws.Range("M139").Value = -2555.2001
$ws.Range("N139").Value = -41420.924

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 15096.9
$ws.Range("I80").Value = 9991
$ws.Range("J80").Value = 17285.143
$ws.Range("K80").Value = 9991
$ws.Range("L80").Value = 17285.143
$ws.Range("M80").Value = -8993
$ws.Range("N80").Value = -19281.143
# Row 83
$ws.Range("H83").Value = 15096.9
$ws.Range("I83").Value = 9991
$ws.Range("J83").Value = 17285.143
$ws.Range("K83").Value = 49955
$ws.Range("L83").Value = 86425.715
$ws.Range("M83").Value = -44963
$ws.Range("N83").Value = -96409.715
# Row 97
$ws.Range("H97").Value = 1009.6667
$ws.Range("J97").Value = 396
$ws.Range("L97").Value = 396
$ws.Range("N97").Value = -1388
# Row 122
$ws.Range("H122").Value = 3987.32
$ws.Range("I122").Value = 3054.5715
$ws.Range("J122").Value = 8884.25
$ws.Range("K122").Value = 9163.7145
$ws.Range("L122").Value = 26652.75
$ws.Range("M122").Value = -6713.7145
$ws.Range("N122").Value = -31552.75
# Row 126
$ws.Range("H126").Value = 5875
$ws.Range("J126").Value = 8000
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -28940

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1854.4166
$ws.Range("I22").Value = 1810.9333
$ws.Range("J22").Value = 1926.8889
$ws.Range("K22").Value = 1810.9333
$ws.Range("L22").Value = 1926.8889
$ws.Range("M22").Value = -1515.9333
$ws.Range("N22").Value = -2516.8889
# Row 27
$ws.Range("H27").Value = 1854.4166
$ws.Range("I27").Value = 1810.9333
$ws.Range("J27").Value = 1926.8889
$ws.Range("K27").Value = 1810.9333
$ws.Range("L27").Value = 1926.8889
$ws.Range("M27").Value = -1703.9333
$ws.Range("N27").Value = -2140.8889
# Row 40
$ws.Range("H40").Value = 9437.409
$ws.Range("I40").Value = 4939.857
$ws.Range("K40").Value = 4939.857
$ws.Range("M40").Value = -4803.857
# Row 46
$ws.Range("H46").Value = 4010.9524
$ws.Range("I46").Value = 4280
$ws.Range("J46").Value = 3845.3845
$ws.Range("K46").Value = 4280
$ws.Range("L46").Value = 3845.3845
$ws.Range("M46").Value = -4092
$ws.Range("N46").Value = -4221.3845
# Row 61
$ws.Range("H61").Value = 4936.6924
$ws.Range("I61").Value = 4931.4165
$ws.Range("K61").Value = 4931.4165
$ws.Range("M61").Value = -4729.4165
# Row 68
$ws.Range("H68").Value = 3077.5
$ws.Range("I68").Value = 1453.2
$ws.Range("K68").Value = 1453.2
$ws.Range("M68").Value = -704.2
# Row 71
$ws.Range("H71").Value = 3077.5
$ws.Range("I71").Value = 1453.2
$ws.Range("K71").Value = 7266
$ws.Range("M71").Value = -3522
# Row 82
$ws.Range("H82").Value = 1352.2222
$ws.Range("J82").Value = 1594.4
$ws.Range("L82").Value = 1594.4
$ws.Range("N82").Value = -2316.4
# Row 85
$ws.Range("H85").Value = 1352.2222
$ws.Range("J85").Value = 1594.4
$ws.Range("L85").Value = 1594.4
$ws.Range("N85").Value = -4090.4
# Row 93
$ws.Range("H93").Value = 2431.2
$ws.Range("I93").Value = 1288.125
$ws.Range("K93").Value = 1288.125
$ws.Range("M93").Value = -40.125
# Row 113
$ws.Range("H113").Value = 4936.6924
$ws.Range("I113").Value = 4931.4165
$ws.Range("K113").Value = 4931.4165
$ws.Range("M113").Value = -2761.4165
# Row 122
$ws.Range("H122").Value = 7205.8
$ws.Range("I122").Value = 7865.4287
$ws.Range("K122").Value = 23596.2861
$ws.Range("M122").Value = -21146.2861
# Row 132
$ws.Range("H132").Value = 4070.138
$ws.Range("I132").Value = 4283.5454
$ws.Range("J132").Value = 3399.4285
$ws.Range("K132").Value = 12850.6362
$ws.Range("L132").Value = 10198.2855
$ws.Range("M132").Value = -10320.6362
$ws.Range("N132").Value = -15258.2855
# Row 136
$ws.Range("H136").Value = 2342.6
$ws.Range("I136").Value = 2013.9
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6041.700000000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3491.700000000001
$ws.Range("N136").Value = -14100

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 26735.111
$ws.Range("J45").Value = 17898.834
$ws.Range("L45").Value = 17898.834
$ws.Range("N45").Value = -18880.834
# Row 62
$ws.Range("H62").Value = 6624
$ws.Range("I62").Value = 5998
$ws.Range("J62").Value = 7250
$ws.Range("K62").Value = 5998
$ws.Range("L62").Value = 7250
$ws.Range("M62").Value = -5374
$ws.Range("N62").Value = -8498
# Row 64
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
# Row 65
$ws.Range("H65").Value = 6624
$ws.Range("I65").Value = 5998
$ws.Range("J65").Value = 7250
$ws.Range("K65").Value = 29990
$ws.Range("L65").Value = 36250
$ws.Range("M65").Value = -26870
$ws.Range("N65").Value = -42490
# Row 67
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
# Row 107
$ws.Range("H107").Value = 3170.4092
$ws.Range("I107").Value = 2667.7144
$ws.Range("J107").Value = 4050.125
$ws.Range("K107").Value = 8003.1432
$ws.Range("L107").Value = 12150.375
$ws.Range("M107").Value = -6083.1432
$ws.Range("N107").Value = -15990.375
# Row 113
$ws.Range("H113").Value = 964.9375
$ws.Range("I113").Value = 964.9375
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2894.8125
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -724.8125
$ws.Range("N113").ClearContents()
# Row 126
$ws.Range("H126").Value = 1941.75
$ws.Range("I126").Value = 1465.1333
$ws.Range("K126").Value = 4395.3999
$ws.Range("M126").Value = -1925.3999
# Row 129
$ws.Range("H129").Value = 96321.75
$ws.Range("J129").Value = 96321.75
$ws.Range("L129").Value = 96321.75
$ws.Range("N129").Value = -106321.75
# Row 132
$ws.Range("H132").Value = 3536.75
$ws.Range("I132").Value = 3215.6667
$ws.Range("K132").Value = 9647.000100000001
$ws.Range("M132").Value = -7117.000100000001
# Row 136
$ws.Range("H136").Value = 1501.44
$ws.Range("I136").Value = 1296.6818
$ws.Range("J136").Value = 3003
$ws.Range("K136").Value = 3890.0454
$ws.Range("L136").Value = 9009
$ws.Range("M136").Value = -1340.0454
$ws.Range("N136").Value = -14109

